$d = $word.ActiveDocument

# --- 1) Merge the two runs that spell "Body<add here your text" + ">" into a single run,
#        dropping the _GoBack bookmark that currently sits between them. ---
$firstRunRange = $d.Content
$null = $firstRunRange.Find.Execute("Body<add here your text", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$mergePoint = $firstRunRange.End

$tailSearch = $d.Range($mergePoint, $d.Content.End)
$null = $tailSearch.Find.Execute(">", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$insertPoint = $d.Range($mergePoint, $mergePoint)
$insertPoint.InsertAfter(">")
$tailSearch.Text = ""

# --- 2) Add a new paragraph with "16:28", and another empty paragraph right after the
#        (now single-run) body paragraph, using Find & Replace so no placeholder run is
#        created for the trailing empty paragraph. ---
$full = $d.Content
$null = $full.Find.Execute("Body<add here your text>", $false, $false, $false, $false, $false, $true, 1, $false, "Body<add here your text>^p16:28^p", 2)

# --- 3) Re-anchor the _GoBack bookmark onto the new, empty trailing paragraph. ---
$lastPara = $d.Paragraphs($d.Paragraphs.Count).Range
$d.Bookmarks.Add("_GoBack", $lastPara)
